# Add a "tbl" column (C) to every sheet that classifies each field as
# belonging to the "main" insert or to a "param" (sub) table, mirroring the
# commit "add a function that select the column names".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 - person
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(1,3).Value = "tbl"
$ws1.Cells.Item(3,3).Value = "main"
$ws1.Cells.Item(4,3).Value = "main"
$ws1.Cells.Item(5,3).Value = "main"
$ws1.Cells.Item(7,3).Value = "main"
$ws1.Cells.Item(14,3).Value = "main"
$ws1.Cells.Item(15,3).Value = "main"
$ws1.Range("C7").Select()

# ---------------------------------------------------------------------
# Sheet 2 - site
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(1,3).Value = "tbl"
$ws2.Cells.Item(3,3).Value = "main"
$ws2.Cells.Item(4,3).Value = "main"
$ws2.Cells.Item(5,3).Value = "main"
$ws2.Cells.Item(6,3).Value = "main"
$ws2.Cells.Item(7,3).Value = "main"
$ws2.Cells.Item(8,3).Value = "main"
$ws2.Cells.Item(9,3).Value = "main"
$ws2.Cells.Item(10,3).Value = "param"
$ws2.Cells.Item(11,3).Value = "main"
$ws2.Cells.Item(12,3).Value = "main"
$ws2.Cells.Item(13,3).Value = "param"
$ws2.Cells.Item(14,3).Value = "param"
$ws2.Cells.Item(15,3).Value = "param"
$ws2.Cells.Item(16,3).Value = "param"
$ws2.Cells.Item(17,3).Value = "param"
$ws2.Cells.Item(18,3).Value = "param"
$ws2.Cells.Item(19,3).Value = "param"
$ws2.Range("C11:C12").Select()

# ---------------------------------------------------------------------
# Sheet 3 - tree
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(1,3).Value = "tbl"
$ws3.Cells.Item(2,3).Value = "main"
$ws3.Cells.Item(4,3).Value = "main"
$ws3.Cells.Item(5,3).Value = "param"
$ws3.Cells.Item(6,3).Value = "param"
$ws3.Cells.Item(7,3).Value = "param"
$ws3.Cells.Item(8,3).Value = "param"
$ws3.Cells.Item(9,3).Value = "param"
$ws3.Cells.Item(10,3).Value = "param"
$ws3.Cells.Item(11,3).Value = "param"
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
$ws3.Range("C6:C11").Select()

# ---------------------------------------------------------------------
# Sheet 4 - sample
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(1,3).Value = "tbl"
$ws4.Cells.Item(3,3).Value = "main"
$ws4.Cells.Item(4,3).Value = "main"
$ws4.Cells.Item(5,3).Value = "param"
$ws4.Cells.Item(6,3).Value = "param"
$ws4.Cells.Item(7,3).Value = "param"
$ws4.Cells.Item(8,3).Value = "param"
$ws4.Cells.Item(9,3).Value = "param"
$ws4.Columns("B").ColumnWidth = 16.5

# ---------------------------------------------------------------------
# Sheet 5 - subsample
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Cells.Item(1,3).Value = "tbl"
$ws5.Cells.Item(4,3).Value = "main"
$ws5.Cells.Item(5,3).Value = "param"
$ws5.Cells.Item(6,3).Value = "param"
$ws5.Cells.Item(7,3).Value = "param"
$ws5.Cells.Item(8,3).Value = "param"
$ws5.Cells.Item(9,3).Value = "param"
$ws5.Cells.Item(10,3).Value = "param"
$ws5.Cells.Item(11,3).Value = "param"
$ws5.Cells.Item(12,3).Value = "param"
$ws5.Cells.Item(13,3).Value = "param"
$ws5.Cells.Item(14,3).Value = "param"
$ws5.Cells.Item(15,3).Value = "param"
$ws5.Cells.Item(16,3).Value = "param"
$ws5.Range("C6:C16").Select()

# ---------------------------------------------------------------------
# Final state: "sample" sheet (4th tab) ends up active with C9 selected
# ---------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("C9").Select()
